$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Column B/C are plain text (coin names / URLs); D/E are numeric-looking
# strings (price / percent change) that must be forced to Text format so
# Excel does not auto-convert them into numbers/percentages.

Set-TextCell $ws 'D2' '256.46'
Set-TextCell $ws 'E2' '4.53%'
Set-TextCell $ws 'D3' '27.58'
Set-TextCell $ws 'E3' '-2.58%'
Set-TextCell $ws 'D4' '5.216'
Set-TextCell $ws 'E4' '-0.70%'
Set-TextCell $ws 'D5' '0.05920'
Set-TextCell $ws 'E5' '3.82%'
Set-TextCell $ws 'D6' '6.690'
Set-TextCell $ws 'E6' '0.89%'
Set-TextCell $ws 'D7' '0.8675'
Set-TextCell $ws 'E7' '1.95%'
Set-TextCell $ws 'D8' '1.030'
Set-TextCell $ws 'E8' '12.40%'
Set-TextCell $ws 'B9' 'WazirX'
Set-TextCell $ws 'C9' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell $ws 'D9' '0.1420'
Set-TextCell $ws 'E9' '3.75%'
Set-TextCell $ws 'B10' 'LiechtensteinCryptoassetsExchange'
Set-TextCell $ws 'C10' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell $ws 'D10' '0.03577'
Set-TextCell $ws 'E10' '8.33%'
Set-TextCell $ws 'D11' '0.07196'
Set-TextCell $ws 'E11' '1.68%'
Set-TextCell $ws 'D12' '0.03263'
Set-TextCell $ws 'E12' '2.43%'
Set-TextCell $ws 'D13' '0.09205'
Set-TextCell $ws 'E13' '-0.17%'
Set-TextCell $ws 'D14' '0.001548'
Set-TextCell $ws 'E14' '1.34%'
Set-TextCell $ws 'B15' 'One'
Set-TextCell $ws 'C15' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextCell $ws 'D15' '0.0006066'
Set-TextCell $ws 'E15' '1.85%'
Set-TextCell $ws 'B16' 'TigerCash'
Set-TextCell $ws 'C16' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell $ws 'D16' '0.005839'
Set-TextCell $ws 'E16' '-1.28%'
Set-TextCell $ws 'B17' 'LEO'
Set-TextCell $ws 'C17' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell $ws 'D17' '3.485'
Set-TextCell $ws 'E17' '-0.16%'
Set-TextCell $ws 'B18' 'GateToken'
Set-TextCell $ws 'C18' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell $ws 'D18' '3.270'
Set-TextCell $ws 'E18' '2.12%'
Set-TextCell $ws 'B19' 'BTSEToken'
Set-TextCell $ws 'C19' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell $ws 'D19' '2.206'
Set-TextCell $ws 'E19' '0.81%'
Set-TextCell $ws 'B20' 'BitpandaEcosystemToken'
Set-TextCell $ws 'C20' 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextCell $ws 'D20' '0.3151'
Set-TextCell $ws 'E20' '-0.60%'
Set-TextCell $ws 'E21' '2.50%'
Set-TextCell $ws 'D22' '3.536'
Set-TextCell $ws 'E22' '0.48%'
Set-TextCell $ws 'D23' '0.04169'
Set-TextCell $ws 'E23' '2.33%'
Set-TextCell $ws 'E24' '1.61%'
Set-TextCell $ws 'D25' '0.001222'
Set-TextCell $ws 'E25' '0.05%'
Set-TextCell $ws 'D26' '0.004518'
Set-TextCell $ws 'E26' '8.70%'
Set-TextCell $ws 'E27' '0.15%'
Set-TextCell $ws 'D28' '0.0001940'
Set-TextCell $ws 'E28' '33.99%'
Set-TextCell $ws 'D40' '0.03817'
Set-TextCell $ws 'E40' '-0.07%'
Set-TextCell $ws 'B41' 'BKEXToken'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextCell $ws 'D41' '0.1104'
Set-TextCell $ws 'E41' '3.57%'
Set-TextCell $ws 'B42' 'KickToken'
Set-TextCell $ws 'C42' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextCell $ws 'D42' '0.003962'
Set-TextCell $ws 'E42' '-23.26%'
Set-TextCell $ws 'D43' '0.002463'
Set-TextCell $ws 'E43' '11.98%'
Set-TextCell $ws 'D44' '0.009840'
Set-TextCell $ws 'E44' '7.41%'
Set-TextCell $ws 'D45' '0.00005437'
Set-TextCell $ws 'E45' '3.35%'
Set-TextCell $ws 'D46' '0.00000000751'
Set-TextCell $ws 'E46' '0.14%'
Set-TextCell $ws 'E47' '4.09%'
Set-TextCell $ws 'E48' '-4.70%'
Set-TextCell $ws 'D49' '0.00002102'
Set-TextCell $ws 'E49' '0.14%'
Set-TextCell $ws 'D50' '0.0002002'
Set-TextCell $ws 'E50' '0.14%'
